{"js": "// Update the division-problem worksheet table: each populated row of the\n// table (rows 0, 4, 8, 12, 16 of the 20-row grid -- the other rows are\n// spacer rows) has its 5 \"NN\u00f7N=\" cells replaced with new problems, in\n// document order, left to right, top to bottom.\n\nconst table = context.document.body.tables.getFirst();\n\n// New values keyed by the table's own (0-based) row/column indices.\nconst updates = [\n  { row: 0, col: 0, text: \"93\u00f79=\" },\n  { row: 0, col: 1, text: \"75\u00f78=\" },\n  { row: 0, col: 2, text: \"49\u00f72=\" },\n  { row: 0, col: 3, text: \"31\u00f72=\" },\n  { row: 0, col: 4, text: \"63\u00f76=\" },\n\n  { row: 4, col: 0, text: \"46\u00f79=\" },\n  { row: 4, col: 1, text: \"92\u00f78=\" },\n  { row: 4, col: 2, text: \"18\u00f78=\" },\n  { row: 4, col: 3, text: \"43\u00f76=\" },\n  { row: 4, col: 4, text: \"94\u00f74=\" },\n\n  { row: 8, col: 0, text: \"58\u00f76=\" },\n  { row: 8, col: 1, text: \"71\u00f78=\" },\n  { row: 8, col: 2, text: \"27\u00f72=\" },\n  { row: 8, col: 3, text: \"20\u00f78=\" },\n  { row: 8, col: 4, text: \"11\u00f77=\" },\n\n  { row: 12, col: 0, text: \"71\u00f74=\" },\n  { row: 12, col: 1, text: \"81\u00f75=\" },\n  { row: 12, col: 2, text: \"89\u00f73=\" },\n  { row: 12, col: 3, text: \"29\u00f77=\" },\n  { row: 12, col: 4, text: \"74\u00f76=\" },\n\n  { row: 16, col: 0, text: \"90\u00f77=\" },\n  { row: 16, col: 1, text: \"72\u00f79=\" },\n  { row: 16, col: 2, text: \"36\u00f77=\" },\n  { row: 16, col: 3, text: \"48\u00f72=\" },\n  { row: 16, col: 4, text: \"37\u00f79=\" },\n];\n\nfor (const { row, col, text } of updates) {\n  table.getCell(row, col).value = text;\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem worksheet table: each populated row of the\n# table (1-based COM rows 1, 5, 9, 13, 17 -- the other rows are spacer rows)\n# has its 5 \"NN\u00f7N=\" cells replaced with new problems, in document order,\n# left to right, top to bottom.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"93\u00f79=\"\n$t.Cell(1, 2).Range.Text = \"75\u00f78=\"\n$t.Cell(1, 3).Range.Text = \"49\u00f72=\"\n$t.Cell(1, 4).Range.Text = \"31\u00f72=\"\n$t.Cell(1, 5).Range.Text = \"63\u00f76=\"\n\n$t.Cell(5, 1).Range.Text = \"46\u00f79=\"\n$t.Cell(5, 2).Range.Text = \"92\u00f78=\"\n$t.Cell(5, 3).Range.Text = \"18\u00f78=\"\n$t.Cell(5, 4).Range.Text = \"43\u00f76=\"\n$t.Cell(5, 5).Range.Text = \"94\u00f74=\"\n\n$t.Cell(9, 1).Range.Text = \"58\u00f76=\"\n$t.Cell(9, 2).Range.Text = \"71\u00f78=\"\n$t.Cell(9, 3).Range.Text = \"27\u00f72=\"\n$t.Cell(9, 4).Range.Text = \"20\u00f78=\"\n$t.Cell(9, 5).Range.Text = \"11\u00f77=\"\n\n$t.Cell(13, 1).Range.Text = \"71\u00f74=\"\n$t.Cell(13, 2).Range.Text = \"81\u00f75=\"\n$t.Cell(13, 3).Range.Text = \"89\u00f73=\"\n$t.Cell(13, 4).Range.Text = \"29\u00f77=\"\n$t.Cell(13, 5).Range.Text = \"74\u00f76=\"\n\n$t.Cell(17, 1).Range.Text = \"90\u00f77=\"\n$t.Cell(17, 2).Range.Text = \"72\u00f79=\"\n$t.Cell(17, 3).Range.Text = \"36\u00f77=\"\n$t.Cell(17, 4).Range.Text = \"48\u00f72=\"\n$t.Cell(17, 5).Range.Text = \"37\u00f79=\"\n"}
